# Weekly data refresh: prepend two new price records (Espárragos, Vega
# Modelo de Temuco) at the top of the data block and push the existing
# records down by two rows. Inserting whole rows preserves every
# existing cell (values + the date-column number format) verbatim in
# their new position, and Excel automatically extends the used range /
# <dimension> to A1:R105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 41:103 down to 43:105, opening up two blank rows at 41:42.
$ws.Range("A41:R42").EntireRow.Insert()

# New row 41
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 45219
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 300000000
$ws.Range("G41").Value = "Espárragos"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Extra"
$ws.Range("J41").Value = 185
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 2000
$ws.Range("M41").Value = 2000
$ws.Range("N41").Value = "$/kilo"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 2000
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"

# New row 42
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 45219
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 300000000
$ws.Range("G42").Value = "Espárragos"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 350
$ws.Range("K42").Value = 1500
$ws.Range("L42").Value = 1500
$ws.Range("M42").Value = 1500
$ws.Range("N42").Value = "$/kilo"
$ws.Range("O42").Value = "Región del Maule"
$ws.Range("P42").Value = 1500
$ws.Range("Q42").Value = 1
$ws.Range("R42").Value = "Hortaliza"
